$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("D1").Value = "NLGDP"

# Single (non-shared) formula in D2
$ws.Range("D2").Formula = "=(LN(B2))"

# Shared-formula groups matching the fill batches used when the column
# was originally created (64-row chunks, last chunk shorter)
$ws.Range("D3:D66").Formula = "=(LN(B3))"
$ws.Range("D67:D130").Formula = "=(LN(B67))"
$ws.Range("D131:D181").Formula = "=(LN(B131))"

# Select column D (mirrors clicking the column header)
[void]$ws.Columns.Item(4).Select()
